# Update LR-pair metrics with recomputed TPM-based values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.597878666666666
$ws.Range("H2").Value = 4.793635999999999
$ws.Range("I2").Value = 0.8992131381376172
$ws.Range("J2").Value = 0.8992131381376171
$ws.Range("M2").Value = 1.484826
$ws.Range("N2").Value = 4.454478
$ws.Range("O2").Value = 0.06049021884829667
$ws.Range("P2").Value = 0.06049021884829667
$ws.Range("Q2").Value = 2.372571789112
$ws.Range("R2").Value = 21.353146102008
$ws.Range("S2").Value = 0.05439359951720808
$ws.Range("T2").Value = 0.05439359951720808

# Row 3
$ws.Range("G3").Value = 1.597878666666666
$ws.Range("H3").Value = 4.793635999999999
$ws.Range("I3").Value = 0.8992131381376172
$ws.Range("J3").Value = 0.8992131381376171
$ws.Range("O3").Value = 0.5859425360316464
$ws.Range("P3").Value = 0.5859425360316464
$ws.Range("S3").Value = 0.5268872265933305
$ws.Range("T3").Value = 0.5268872265933304

# Row 4
$ws.Range("G4").Value = 1.597878666666666
$ws.Range("H4").Value = 4.793635999999999
$ws.Range("I4").Value = 0.8992131381376172
$ws.Range("J4").Value = 0.8992131381376171
$ws.Range("M4").Value = 8.653369666666666
$ws.Range("N4").Value = 25.960109
$ws.Range("O4").Value = 0.3525289999716321
$ws.Range("P4").Value = 0.3525289999716321
$ws.Range("Q4").Value = 13.82703478514711
$ws.Range("R4").Value = 124.443313066324
$ws.Range("S4").Value = 0.3169987083490073
$ws.Range("T4").Value = 0.3169987083490072

# Row 5
$ws.Range("G5").Value = 1.597878666666666
$ws.Range("H5").Value = 4.793635999999999
$ws.Range("I5").Value = 0.8992131381376172
$ws.Range("J5").Value = 0.8992131381376171
$ws.Range("M5").Value = 0.02548533333333333
$ws.Range("N5").Value = 0.076456
$ws.Range("O5").Value = 0.001038245148424882
$ws.Range("P5").Value = 0.001038245148424882
$ws.Range("Q5").Value = 0.04072247044622222
$ws.Range("R5").Value = 0.3665022340159999
$ws.Range("S5").Value = 0.000933603678071294
$ws.Range("T5").Value = 0.0009336036780712939

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.1790956666666667
$ws.Range("H6").Value = 0.537287
$ws.Range("I6").Value = 0.1007868618623829
$ws.Range("J6").Value = 0.1007868618623829
$ws.Range("M6").Value = 1.484826
$ws.Range("N6").Value = 4.454478
$ws.Range("O6").Value = 0.06049021884829667
$ws.Range("P6").Value = 0.06049021884829667
$ws.Range("Q6").Value = 0.265925902354
$ws.Range("R6").Value = 2.393333121186
$ws.Range("S6").Value = 0.00609661933108859
$ws.Range("T6").Value = 0.006096619331088589

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.1790956666666667
$ws.Range("H7").Value = 0.537287
$ws.Range("I7").Value = 0.1007868618623829
$ws.Range("J7").Value = 0.1007868618623829
$ws.Range("O7").Value = 0.5859425360316464
$ws.Range("P7").Value = 0.5859425360316464
$ws.Range("Q7").Value = 2.575908974847333
$ws.Range("R7").Value = 23.183180773626
$ws.Range("S7").Value = 0.05905530943831588
$ws.Range("T7").Value = 0.05905530943831588

# Row 8
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.1790956666666667
$ws.Range("H8").Value = 0.537287
$ws.Range("I8").Value = 0.1007868618623829
$ws.Range("J8").Value = 0.1007868618623829
$ws.Range("M8").Value = 8.653369666666666
$ws.Range("N8").Value = 25.960109
$ws.Range("O8").Value = 0.3525289999716321
$ws.Range("P8").Value = 0.3525289999716321
$ws.Range("Q8").Value = 1.549781009364778
$ws.Range("R8").Value = 13.948029084283
$ws.Range("S8").Value = 0.03553029162262489
$ws.Range("T8").Value = 0.03553029162262488

# Row 9
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.1790956666666667
$ws.Range("H9").Value = 0.537287
$ws.Range("I9").Value = 0.1007868618623829
$ws.Range("J9").Value = 0.1007868618623829
$ws.Range("M9").Value = 0.02548533333333333
$ws.Range("N9").Value = 0.076456
$ws.Range("O9").Value = 0.001038245148424882
$ws.Range("P9").Value = 0.001038245148424882
$ws.Range("Q9").Value = 0.004564312763555555
$ws.Range("R9").Value = 0.041078814872
$ws.Range("S9").Value = 0.0001046414703535878
$ws.Range("T9").Value = 0.0001046414703535878

